$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.994.12'
$ws.Range("E2").Value = '  +0.40%  '

$ws.Range("D3").Value = '1.560.88'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '

$ws.Range("E6").Value = '  +0.92%  '

$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.14'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.21%  '

$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0597'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.02%  '

$ws.Range("E11").Value = '  +0.24%  '

$ws.Range("D12").Value = '1.783.14'
$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("D13").Value = '1.541.81'
$ws.Range("E13").Value = '  -0.67%  '

$ws.Range("E14").Value = '  +1.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.78%  '

$ws.Range("D17").Value = '26.976.39'
$ws.Range("E17").Value = '  +0.30%  '

$ws.Range("D18").Value = '0.0₃0706'
$ws.Range("E18").Value = '  +2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.73%  '

$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("E24").Value = '  -1.58%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("E28").Value = '  +1.45%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("E31").Value = '  +1.91%  '

$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").Value = '1.426.08'
$ws.Range("E33").Value = '  +0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.07%  '

$ws.Range("E35").Value = '  +3.28%  '

$ws.Range("E36").Value = '  +9.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0166'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.531'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.809'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.37%  '

$ws.Range("E41").Value = '  -0.17%  '

$ws.Range("E42").Value = '  -0.29%  '

$ws.Range("E43").Value = '  +2.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.14%  '

$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("D47").Value = '1.696.59'
$ws.Range("E47").Value = '  +0.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.44%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  -1.15%  '

$ws.Range("E51").Value = '  -0.02%  '
